$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.397975444793701
$ws.Range("B1").Value = 2.67406702041626
$ws.Range("C1").Value = 2.871906757354736
$ws.Range("D1").Value = 3.174413204193115
$ws.Range("E1").Value = 0.8022944927215576
